$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for Price/Volume columns so values are not
# reinterpreted as numbers or percentages by Excel.
$rng = $ws.Range("D2:E51")
$rng.NumberFormat = "@"

$ws.Range("D2").Value = "328.00"
$ws.Range("E2").Value = "-0.84%"

$ws.Range("D3").Value = "43.86"
$ws.Range("E3").Value = "5.41%"

$ws.Range("D4").Value = "5.413"
$ws.Range("E4").Value = "-4.86%"

$ws.Range("D5").Value = "0.08091"
$ws.Range("E5").Value = "-3.55%"

$ws.Range("D6").Value = "8.703"
$ws.Range("E6").Value = "-1.18%"

$ws.Range("D7").Value = "1.913"
$ws.Range("E7").Value = "-4.79%"

$ws.Range("D8").Value = "4.303"
$ws.Range("E8").Value = "-3.65%"

$ws.Range("E9").Value = "-4.75%"

$ws.Range("D10").Value = "0.9423"
$ws.Range("E10").Value = "1.97%"

$ws.Range("D11").Value = "0.1206"
$ws.Range("E11").Value = "-5.93%"

$ws.Range("D12").Value = "0.1899"
$ws.Range("E12").Value = "-3.89%"

$ws.Range("D13").Value = "0.09520"
$ws.Range("E13").Value = "-0.27%"

$ws.Range("D14").Value = "0.04171"
$ws.Range("E14").Value = "8.53%"

$ws.Range("D15").Value = "0.1070"
$ws.Range("E15").Value = "0.68%"

$ws.Range("E16").Value = "-1.41%"

$ws.Range("D17").Value = "0.006064"
$ws.Range("E17").Value = "-0.64%"

$ws.Range("D18").Value = "3.570"
$ws.Range("E18").Value = "4.20%"

$ws.Range("D20").Value = "8.518"
$ws.Range("E20").Value = "-3.58%"

$ws.Range("D21").Value = "0.1358"
$ws.Range("E21").Value = "-0.35%"

$ws.Range("E22").Value = "3.83%"

$ws.Range("D23").Value = "0.04378"
$ws.Range("E23").Value = "-0.48%"

$ws.Range("E24").Value = "-2.51%"

$ws.Range("D25").Value = "0.004314"
$ws.Range("E25").Value = "-1.47%"

$ws.Range("D26").Value = "0.0001239"
$ws.Range("E26").Value = "1.54%"

$ws.Range("D27").Value = "0.0004016"
$ws.Range("E27").Value = "0.68%"

$ws.Range("D39").Value = "0.02675"
$ws.Range("E39").Value = "-6.93%"

$ws.Range("D40").Value = "0.05455"
$ws.Range("E40").Value = "-1.29%"

$ws.Range("D41").Value = "0.007800"
$ws.Range("E41").Value = "-1.97%"

$ws.Range("D42").Value = "0.009788"
$ws.Range("E42").Value = "8.75%"

$ws.Range("D43").Value = "0.1394"
$ws.Range("E43").Value = "-2.88%"

$ws.Range("D44").Value = "0.002126"
$ws.Range("E44").Value = "2.73%"

$ws.Range("D45").Value = "0.009965"
$ws.Range("E45").Value = "-14.62%"

$ws.Range("D46").Value = "0.00007078"
$ws.Range("E46").Value = "2.09%"

$ws.Range("E47").Value = "0.69%"

$ws.Range("D48").Value = "0.003473"
$ws.Range("E48").Value = "0.26%"

$ws.Range("D49").Value = "0.002285"
$ws.Range("E49").Value = "0.35%"

$ws.Range("D50").Value = "0.00002114"
$ws.Range("E50").Value = "0.69%"

$ws.Range("D51").Value = "0.0002014"
$ws.Range("E51").Value = "0.69%"

# Restore default (General) formatting so no stray styles remain.
$rng.ClearFormats()
